$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: RSQUO (right single quotation mark U+2019) used throughout.
# ---------------------------------------------------------------------
$rsquo = [char]0x2019

# ---------------------------------------------------------------------
# Change 1: add <w:rPr><w:noProof/></w:rPr> to the run that holds the
# screenshot drawing (6th paragraph).
# ---------------------------------------------------------------------
$picPara = $d.Paragraphs.Item(6)
$picShape = $picPara.Range.InlineShapes.Item(1)
$picShape.Range.Font.NoProofing = 1

# ---------------------------------------------------------------------
# A reference range carrying "Bold + BoldCS" (w:b + w:bCs) formatting,
# taken from the existing "Send" run, so we can stamp the same pair of
# properties onto newly-created bold runs ("Re-upload" / "Notes") further
# down. (Directly setting .Font.Bold only ever emits <w:b/>, never the
# accompanying <w:bCs/>, so we copy formatting from an existing run that
# already carries both.)
# ---------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("Send to EDI", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$boldSrcRange = $d.Range($findRng.Start, $findRng.Start + 4)

function Set-BoldLikeSend($rng) {
    $rng.FormattedText = $boldSrcRange.FormattedText
}

# ---------------------------------------------------------------------
# Change 2: the paragraph that used to read
#   "Enter your name and email address. If you have any notes or
#    questions for EDI's data curation team, enter them in the Notes
#    text area."
# becomes three paragraphs:
#   1) a new paragraph about re-uploading data files (bold "Re-upload")
#   2) a blank paragraph
#   3) the reworded "Enter your name..." paragraph (bold "Notes")
# ---------------------------------------------------------------------
$enterPara = $d.Paragraphs.Item(8)

# Make room: insert two blank paragraphs before paragraph 8. Afterwards:
#   paragraph 8  -> new, empty (becomes the "This page reminds..." text)
#   paragraph 9  -> new, empty (stays blank)
#   paragraph 10 -> original "Enter your name..." paragraph (to be reworded)
$d.Range($enterPara.Range.Start, $enterPara.Range.Start).InsertParagraphBefore() | Out-Null
$d.Range($enterPara.Range.Start, $enterPara.Range.Start).InsertParagraphBefore() | Out-Null

# --- fill paragraph 8 ("This page reminds you ... Re-upload ... out-of-date.") ---
$p8 = $d.Paragraphs.Item(8)
$pos = $p8.Range.Start

$t1 = "This page reminds you that the data files (data tables and other entities) that you have uploaded are included in the package sent to EDI, so you will want to make sure you have uploaded the latest versions of all such files. If necessary, use the "
$d.Range($pos, $pos).InsertAfter($t1) | Out-Null
$pos = $pos + $t1.Length

$t2 = "Re-upload"
$d.Range($pos, $pos).InsertAfter($t2) | Out-Null
$boldRng = $d.Range($pos, $pos + $t2.Length)
Set-BoldLikeSend $boldRng
$fixRng = $d.Range($pos, $pos + $t2.Length)
$fixRng.Text = $t2
$pos = $pos + $t2.Length

$t3 = " features to update any files that are out-of-date."
$d.Range($pos, $pos).InsertAfter($t3) | Out-Null
$pos = $pos + $t3.Length

# --- paragraph 9 stays blank (already inserted as empty paragraph) ---

# --- reword paragraph 10 ("Enter your name ... Notes ... text area.") ---
$p10 = $d.Paragraphs.Item(10)
$start = $p10.Range.Start
$end = $p10.Range.End
$d.Range($start, $end - 1).Text = ""

$pos = $start
$u1 = "Enter your name and email address. If you have any notes or questions for EDI" + $rsquo + "s data curation team, enter them in the "
$d.Range($pos, $pos).InsertAfter($u1) | Out-Null
$pos = $pos + $u1.Length

$u2 = "Notes"
$d.Range($pos, $pos).InsertAfter($u2) | Out-Null
$boldRng2 = $d.Range($pos, $pos + $u2.Length)
Set-BoldLikeSend $boldRng2
$fixRng2 = $d.Range($pos, $pos + $u2.Length)
$fixRng2.Text = $u2
$pos = $pos + $u2.Length

$u3 = " text area."
$d.Range($pos, $pos).InsertAfter($u3) | Out-Null
$pos = $pos + $u3.Length

# ---------------------------------------------------------------------
# Change 3: the paragraph
#   "Please note that there will be some amount of delay before your
#    data package appears in EDI's Data Portal."
# becomes the same sentence, reworded, but typed/assembled as five
# separate (same-formatting) runs:
#   "Please note that" / " because of these manual curation steps, the
#   data" / " " / "package will not show up immediately " /
#   "in EDI's Data Portal."
# A temporary bookmark is dropped after each inserted chunk and removed
# right away; this forces Word to keep the runs distinct instead of
# silently re-coalescing adjacent same-format runs when the file is
# saved.
# ---------------------------------------------------------------------
$noticePara = $d.Paragraphs.Item(14)
$start = $noticePara.Range.Start
$end = $noticePara.Range.End
$d.Range($start, $end - 1).Text = ""

$parts = @(
    "Please note that",
    " because of these manual curation steps, the data",
    " ",
    "package will not show up immediately ",
    ("in EDI" + $rsquo + "s Data Portal.")
)

$pos = $start
$bmCount = 0
foreach ($part in $parts) {
    $d.Range($pos, $pos).InsertAfter($part) | Out-Null
    $newPos = $pos + $part.Length
    $bmCount = $bmCount + 1
    $d.Bookmarks.Add("tmpRunBreak" + $bmCount, $d.Range($pos, $newPos)) | Out-Null
    $pos = $newPos
}
for ($i = 1; $i -le $bmCount; $i++) {
    $d.Bookmarks.Item("tmpRunBreak" + $i).Delete()
}
